# Update the "Förändrad" (Changed) date column (C) from 2023-10-08 (45207)
# to 2023-10-09 (45208) for all data rows (rows 2 through 158).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 158

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2() -eq 45207) {
        $cell.Value = 45208
    }
}
